$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before row 2 (the current Nokron row), shifting
# the existing Nokron/Siofra data down from rows 2-6 to rows 5-9.
$ws.Range("A2:C4").Insert()
# The insert picks up formatting from the row above (the bold header);
# the newly-inserted data rows should be unformatted, like the rest of
# the data rows, so strip that back off.
$ws.Range("A2:C4").ClearFormats()

# New rows: Ainsel River data
$ws.Range("A2").Value = "Ainsel River"
$ws.Range("B2").Value = "Major bosses"
$ws.Range("C2").Value = "Dropped by Astel, Naturalborn of the Void. Replaces Remembrance of the Naturalborn"

$ws.Range("A3").Value = "Ainsel River Downstream"
$ws.Range("B3").Value = "Merchant shops"
$ws.Range("C3").Value = "Sold by the Hermit Merchant in the alcove past the rock-slinging Malformed Star"

$ws.Range("A4").Value = "Ainsel River Downstream"
$ws.Range("B4").Value = "Major bosses"
$ws.Range("C4").Value = "Dropped by Dragonkin Soldier of Nokstella. Replaces Frozen Lightning Spear"
